# EPBDS-8724: Simple Rules and Simple Lookups should not support StringRanges
# for backward compatibility.
#
# The "SimpleRules"/"Test mySimpeRulesN" example tables (rows 36-39 and
# 46-49 on Sheet1) used string-range values such as "AAA - BBB" which are
# no longer a supported scenario for Simple Rules/Simple Lookups, so the
# sample data + headers for those tables are wiped out (content, formatting
# and the merged header cells), leaving the now-empty rows behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First table block ("SimpleRules ... mySimpeRulesN" headers/sample data).
$ws.Range("B36:J39").UnMerge()
$ws.Range("B36:J39").Clear()

# Second table block ("Test mySimpeRulesN" headers/sample data).
$ws.Range("B46:J49").UnMerge()
$ws.Range("B46:J49").Clear()
